$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows right after the current last data row (21). Inserting
# (rather than just writing into empty rows below the used range) makes
# Excel carry the existing row's per-column number formats down onto the
# new rows, so the new data lines up with the dd/mm/yyyy date column and
# the '@' text columns used throughout the rest of the table.
$ws.Rows("22:24").Insert()

$dt = Get-Date -Year 2018 -Month 5 -Day 15 -Hour 0 -Minute 0 -Second 0 -Millisecond 0

$ws.Range("A22").Value = $dt
$ws.Range("B22").Value = "Saldakeeva Elena"
$ws.Range("C22").Value = "Rocche filo nero"
$ws.Range("D22").Value = "N°."
$ws.Range("E22").Value = 7

$ws.Range("A23").Value = $dt
$ws.Range("B23").Value = "Saldakeeva Elena"
$ws.Range("C23").Value = "Rocche filo bianco "
$ws.Range("D23").Value = "N°."
$ws.Range("E23").Value = 10

$ws.Range("A24").Value = $dt
$ws.Range("B24").Value = "Segreteria"
$ws.Range("C24").Value = "Squadretta"
$ws.Range("D24").Value = "N°."
$ws.Range("E24").Value = 5
